$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '56.760.09'
$ws.Range('E2').Value = '  -0.15%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.412.77'
$ws.Range('E3').Value = '  -3.61%  '

# Row 4
$ws.Range('E4').Value = '  +0.19%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '486.21'
$ws.Range('E5').Value = '  -2.00%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '153.13'
$ws.Range('E6').Value = '  -0.27%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.12%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.611'
$ws.Range('E8').Value = '  +18.19%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.430.38'
$ws.Range('E9').Value = '  -3.34%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0996'
$ws.Range('E10').Value = '  +0.31%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.64'
$ws.Range('E11').Value = '  -2.42%  '

# Row 12
$ws.Range('E12').Value = '  -0.76%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.126'
$ws.Range('E13').Value = '  +1.03%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.842.28'

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '56.992.65'
$ws.Range('E15').Value = '  +0.28%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.73'
$ws.Range('E16').Value = '  -3.32%  '

# Row 17
$ws.Range('E17').Value = '  -3.09%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.434.02'
$ws.Range('E18').Value = '  -3.20%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.73'
$ws.Range('E19').Value = '  +3.55%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '324.60'
$ws.Range('E20').Value = '  +0.03%  '

# Row 21
$ws.Range('E21').Value = '  -3.73%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  +0.16%  '

# Row 23
$ws.Range('E23').Value = '  -0.92%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '58.14'
$ws.Range('E24').Value = '  -1.62%  '

# Row 25
$ws.Range('E25').Value = '  -1.04%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.42%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.158'
$ws.Range('E27').Value = '  -3.50%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.526.43'
$ws.Range('E28').Value = '  -3.16%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.22'
$ws.Range('E29').Value = '  -6.47%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0782'
$ws.Range('E30').Value = '  -4.36%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  -0.06%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '18.60'
$ws.Range('E32').Value = '  +0.71%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '148.55'
$ws.Range('E33').Value = '  -1.87%  '

# Row 34
$ws.Range('E34').Value = '  -0.87%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.30'
$ws.Range('E35').Value = '  +0.64%  '

# Row 36
$ws.Range('E36').Value = '  -2.37%  '

# Row 37
$ws.Range('E37').Value = '  -3.04%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.847'
$ws.Range('E38').Value = '  -4.18%  '

# Row 39
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '34.12'
$ws.Range('E39').Value = '  -0.33%  '

# Row 40
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.101'
$ws.Range('E40').Value = '  +9.10%  '

# Row 41
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.52'
$ws.Range('E41').Value = '  -0.25%  '

# Row 42
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.37'
$ws.Range('E42').Value = '  -2.28%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.995'
$ws.Range('E43').Value = '  +0.12%  '

# Row 44
$ws.Range('E44').Value = '  -3.80%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '265.16'
$ws.Range('E45').Value = '  -2.35%  '

# Row 46
$ws.Range('E46').Value = '  -6.27%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.73'
$ws.Range('E47').Value = '  -4.54%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '10.21'
$ws.Range('E48').Value = '  +0.06%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0228'
$ws.Range('E49').Value = '  -1.67%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '17.49'
$ws.Range('E50').Value = '  -2.88%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.865.59'
$ws.Range('E51').Value = '  -2.50%  '
